# Add team record (Wins/Losses/Ties) columns to the roster sheet.
# New columns land right after the existing "Unnamed: 28" column (AC),
# extending the used range from A1:AC57 to A1:AF57.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column numbers: AC=29, AD=30, AE=31, AF=32
$colWins   = 30
$colLosses = 31
$colTies   = 32

# Clone the existing header formatting (bold, centered, bordered) from the
# last header cell (AC1) onto the three new header cells so they match the
# rest of row 1.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item(1, $colWins).Value   = "Wins"
$ws.Cells.Item(1, $colLosses).Value = "Losses"
$ws.Cells.Item(1, $colTies).Value   = "Ties"

# Every player row gets the same team record: 55 wins, 107 losses, 0 ties.
$lastRow = 57
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, $colWins).Value   = 55
    $ws.Cells.Item($r, $colLosses).Value = 107
    $ws.Cells.Item($r, $colTies).Value   = 0
}
